$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 40 that duplicates row 39 (same date + gold price text),
# copying both values/shared-string reuse and formatting (styles 3/4).
$ws.Range("A39:B39").Copy()
$ws.Range("A40").PasteSpecial()
